# StorageComponentClassDiagram.pptx update
#
# Commit: "update Design components in DevelopGuide.adoc"
#
# The only content-level change on the slide is the rename of the
# "XmlAdaptedPerson" class box to "XmlAdaptedParcel" in the storage
# component class diagram (slide 1). Locate the shape by its current
# text (rather than a hard-coded shape index) so the edit is resilient
# to shape ordering, then update just its text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$oldName = "XmlAdaptedPerson"
$newName = "XmlAdaptedParcel"

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $textRange = $shape.TextFrame.TextRange
        if ($textRange.Text -eq $oldName) {
            $textRange.Text = $newName
        }
    }
}
